# lifecycle-to-plant.xlsx — add Cactus (perennial) row + reference note/link

$wb = $excel.ActiveWorkbook
$kb = $wb.Worksheets.Item("KB")
$refs = $wb.Worksheets.Item("References")

# --- KB sheet: append row 117 (Cactus / perennial) ---
$kb.Range("A117").Value = "Cactus"
$kb.Range("B117").Value = "perennial"

# --- References sheet: append explanatory note + hyperlink for row 117 ---
$refs.Range("B24").Value = "For row 117 (Cactus)"

$refs.Range("B25").Value = "https://www.hunker.com/13427982/are-succulents-perennials"
$refs.Hyperlinks.Add($refs.Range("B25"), "https://www.hunker.com/13427982/are-succulents-perennials") | Out-Null
$refs.Range("B25").Style = "Hyperlink"

# --- Selection / view state to match the author's final workbook ---
# Scroll the KB sheet so row 103 is the top visible row, then land the
# selection on D118 (first empty row under the new Cactus entry).
$kb.Activate()
$kb.Application.Goto($kb.Range("A103"), $true)
$kb.Application.ActiveWindow.ScrollRow = 103
$kb.Range("D118").Select()

# References sheet keeps its own pending selection one blank row below
# the newly-added hyperlink row, but KB (tab 1) is the active tab now.
$refs.Range("B27").Select()

$kb.Activate()
